$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B4 values
$ws.Range("B2").Value = 0.18536461478977176
$ws.Range("B3").Value = -0.4718430458799713
$ws.Range("B4").Value = -0.06808545292455209

# Remove row 5 entirely (A5/B5 had "4" and -0.10820192437641177)
$ws.Range("A5:B5").Delete()
